# "Implemented New Design into Page Pagination System + OCR Planning"
#
# The sheet previously used Excel's "Place in Cell" / linked-picture rich
# values (F2:F4 held #VALUE! placeholders backed by rdRichValue /
# richValueRel metadata while the images finished resolving). The new
# design drops that OCR/image-linking experiment for now and replaces
# those placeholder cells with a plain "-" text value, and tidies up the
# row heights that had been auto-sized around the (now removed) inline
# images.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Items")

# Replace the broken rich-value / linked-image placeholders in F2:F4 with
# a simple "-" text marker.
$ws.Range("F2").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "-"

# The big custom row heights were sized for the in-cell images that used
# to live behind those #VALUE! placeholders; auto-fit them back down now
# that the cells just hold short text.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# Move the active selection to reflect where the editor left off.
[void]$ws.Range("F9").Select()

$wb.Save()
